$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.022.45'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = '  -0.31%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.872.38'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = '  -2.46%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '319.55'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = '  -3.17%  '

# Row 6
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5041'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = '  -3.42%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3954'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = '  -3.29%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08202'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = '  -3.79%  '

# Row 10
$ws.Range("E10").Value = '  -2.63%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.092'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = '  -3.22%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '23.68'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = '  +5.61%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.866.72'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = '  -2.99%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.297'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = '  -1.85%  '

# Row 15
$ws.Range("E15").Value = '  -3.24%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.002'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '

# Row 17
$ws.Range("E17").Value = '  -4.13%  '

# Row 18
$ws.Range("E18").Value = '  -2.31%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06409'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = '  -4.60%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.13'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = '  -1.03%  '

# Row 21
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '30.019.69'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = '  -0.38%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.849'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = '  -3.15%  '

# Row 24
$ws.Range("E24").Value = '  -2.05%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.168'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = '  -2.48%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.086.48'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = '  -2.71%  '

# Row 27
$ws.Range("E27").Value = '  +0.91%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '160.31'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = '  +0.16%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.223'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = '  -9.62%  '

# Row 30
$ws.Range("E30").Value = '  -1.48%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.069'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = '  -1.05%  '

# Row 32
$ws.Range("E32").Value = '  -2.23%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.937'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = '  -2.63%  '

# Row 34
$ws.Range("E34").Value = '  +1.17%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.02437'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = '  -2.72%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.210'
$ws.Cells.Item(36, 4).Style = "Normal"

# Row 37
$ws.Range("E37").Value = '  -3.80%  '

# Row 38
$ws.Range("E38").Value = '  -3.26%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.174'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = '  -4.96%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '8.483'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = '  -5.09%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6308'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = '  -3.61%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.218'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = '  -2.46%  '

# Row 43
$ws.Range("E43").Value = '  -3.74%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.9997'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = '  -0.07%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.5911'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = '  -4.32%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.94'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = '  -2.94%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.095'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = '  +0.18%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.627'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = '  -3.94%  '

# Row 49
$ws.Range("E49").Value = '  -1.77%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.205'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = '  -3.72%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '77.54'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = '  -2.87%  '
